$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.46388812271717
$ws.Range("D2").Value = 9.053051090061398
$ws.Range("E2").Value = 13.49814074132279
$ws.Range("F2").Value = 34.16091268773715
$ws.Range("G2").Value = 3.623242451708317
$ws.Range("I2").Value = 42.00370031271751
$ws.Range("J2").Value = 9.848448925951388
$ws.Range("O2").Value = 25.54137488967322
$ws.Range("B3").Value = 14.62529716494162
$ws.Range("D3").Value = 9.062502353047922
$ws.Range("E3").Value = 13.47515106134687
$ws.Range("F3").Value = 33.77373060893175
$ws.Range("G3").Value = 3.626715730940678
$ws.Range("I3").Value = 39.58439637612544
$ws.Range("J3").Value = 9.851853114224063
$ws.Range("O3").Value = 25.2832732662348
$ws.Range("B4").Value = 14.08227189835594
$ws.Range("D4").Value = 9.06961810779516
$ws.Range("E4").Value = 13.46397592245685
$ws.Range("F4").Value = 33.54394766294402
$ws.Range("G4").Value = 3.628960586445301
$ws.Range("I4").Value = 38.02187358361297
$ws.Range("J4").Value = 9.855753084105856
$ws.Range("O4").Value = 25.13087928067921
$ws.Range("B5").Value = 13.85402509196118
$ws.Range("D5").Value = 9.072847943946
$ws.Range("E5").Value = 13.46016374882351
$ws.Range("F5").Value = 33.45241244229403
$ws.Range("G5").Value = 3.629903714509761
$ws.Range("I5").Value = 37.36626242888532
$ws.Range("J5").Value = 9.857797311550719
$ws.Range("O5").Value = 25.07037136719887
$ws.Range("B6").Value = 13.81570835831321
$ws.Range("D6").Value = 9.073404193405828
$ws.Range("E6").Value = 13.45957559886189
$ws.Range("F6").Value = 33.43734307526277
$ws.Range("G6").Value = 3.630062034442656
$ws.Range("I6").Value = 37.25627603765865
$ws.Range("J6").Value = 9.85816422533258
$ws.Range("O6").Value = 25.06042215919854
$ws.Range("B7").Value = 14.07922169308636
$ws.Range("D7").Value = 9.069660329893301
$ws.Range("E7").Value = 13.46392150404266
$ws.Range("F7").Value = 33.54270453687643
$ws.Range("G7").Value = 3.628973190948788
$ws.Range("I7").Value = 38.01310743513349
$ws.Range("J7").Value = 9.855778811462464
$ws.Range("O7").Value = 25.13005671365274
$ws.Range("B8").Value = 15.18063999774328
$ws.Range("D8").Value = 9.056037481920358
$ws.Range("E8").Value = 13.48960515884906
$ws.Range("F8").Value = 34.02582672660026
$ws.Range("G8").Value = 3.624416804934428
$ws.Range("I8").Value = 41.18578900812587
$ws.Range("J8").Value = 9.849247038929464
$ws.Range("O8").Value = 25.45115907982008
$ws.Range("B9").Value = 17.11354650569201
$ws.Range("D9").Value = 9.039735604397224
$ws.Range("E9").Value = 13.56317490237894
$ws.Range("F9").Value = 35.03149178840476
$ws.Range("G9").Value = 3.616367576086778
$ws.Range("I9").Value = 46.77879708202151
$ws.Range("J9").Value = 9.850798242082389
$ws.Range("O9").Value = 26.1260939801174
$ws.Range("B10").Value = 18.39153772040988
$ws.Range("D10").Value = 9.034100795024692
$ws.Range("E10").Value = 13.63116400232306
$ws.Range("F10").Value = 35.79921196626681
$ws.Range("G10").Value = 3.610987116409541
$ws.Range("I10").Value = 50.48748964831373
$ws.Range("J10").Value = 9.86068428796484
$ws.Range("O10").Value = 26.64535503308267
$ws.Range("B11").Value = 18.94157951104755
$ws.Range("D11").Value = 9.032912160003942
$ws.Range("E11").Value = 13.66506339699835
$ws.Range("F11").Value = 36.15327824047178
$ws.Range("G11").Value = 3.608653775659024
$ws.Range("I11").Value = 52.08534703260537
$ws.Range("J11").Value = 9.867075651012993
$ws.Range("O11").Value = 26.88572683328339
$ws.Range("B12").Value = 19.14533481163442
$ws.Range("D12").Value = 9.032659395608954
$ws.Range("E12").Value = 13.67832158929037
$ws.Range("F12").Value = 36.28792533965255
$ws.Range("G12").Value = 3.607786520170254
$ws.Range("I12").Value = 52.67744494454099
$ws.Range("J12").Value = 9.869767470888318
$ws.Range("O12").Value = 26.97726771910306
$ws.Range("B13").Value = 19.10165446781605
$ws.Range("D13").Value = 9.032705061912115
$ws.Range("E13").Value = 13.67544756740442
$ws.Range("F13").Value = 36.25890330013022
$ws.Range("G13").Value = 3.607972574463442
$ws.Range("I13").Value = 52.55050504930303
$ws.Range("J13").Value = 9.869175677268064
$ws.Range("O13").Value = 26.95753102762593
$ws.Range("B14").Value = 18.95843357864365
$ws.Range("D14").Value = 9.032887411950966
$ws.Range("E14").Value = 13.66614575069779
$ws.Range("F14").Value = 36.16434502539727
$ws.Range("G14").Value = 3.608582099321983
$ws.Range("I14").Value = 52.13432002565632
$ws.Range("J14").Value = 9.867291671810763
$ws.Range("O14").Value = 26.89324805125582
$ws.Range("B15").Value = 18.87011552320599
$ws.Range("D15").Value = 9.033024795924055
$ws.Range("E15").Value = 13.66050278847509
$ws.Range("F15").Value = 36.10649589555697
$ws.Range("G15").Value = 3.608957574585209
$ws.Range("I15").Value = 51.87770123785431
$ws.Range("J15").Value = 9.866173000034886
$ws.Range("O15").Value = 26.85393784757937
$ws.Range("B16").Value = 18.3549581758265
$ws.Range("D16").Value = 9.034206110330461
$ws.Range("E16").Value = 13.62900788030495
$ws.Range("F16").Value = 35.7761593962845
$ws.Range("G16").Value = 3.611141896300692
$ws.Range("I16").Value = 50.38125698588405
$ws.Range("J16").Value = 9.8603046461317
$ws.Range("O16").Value = 26.62972294951099
$ws.Range("B17").Value = 18.03087938850338
$ws.Range("D17").Value = 9.035282680972486
$ws.Range("E17").Value = 13.61044360731533
$ws.Range("F17").Value = 35.57465633020601
$ws.Range("G17").Value = 3.612511099626303
$ws.Range("I17").Value = 49.44026495104949
$ws.Range("J17").Value = 9.857189210147087
$ws.Range("O17").Value = 26.49318197498308
$ws.Range("B18").Value = 17.84153069387959
$ws.Range("D18").Value = 9.036031310258103
$ws.Range("E18").Value = 13.60004595055646
$ws.Range("F18").Value = 35.45921881065992
$ws.Range("G18").Value = 3.613309389561931
$ws.Range("I18").Value = 48.89063678530206
$ws.Range("J18").Value = 9.855575641806848
$ws.Range("O18").Value = 26.41504338001178
$ws.Range("B19").Value = 17.77691535278072
$ws.Range("D19").Value = 9.036307019224486
$ws.Range("E19").Value = 13.59657374906437
$ws.Range("F19").Value = 35.42021674247732
$ws.Range("G19").Value = 3.613581527980174
$ws.Range("I19").Value = 48.7031053676155
$ws.Range("J19").Value = 9.855059968827746
$ws.Range("O19").Value = 26.3886574843856
$ws.Range("B20").Value = 18.06568338928438
$ws.Range("D20").Value = 9.035154686480276
$ws.Range("E20").Value = 13.61239086607596
$ws.Range("F20").Value = 35.59605979171533
$ws.Range("G20").Value = 3.61236423266498
$ws.Range("I20").Value = 49.54130510666898
$ws.Range("J20").Value = 9.857502400698968
$ws.Range("O20").Value = 26.50767657635921
$ws.Range("B21").Value = 19.00062425561376
$ws.Range("D21").Value = 9.032828498380381
$ws.Range("E21").Value = 13.66886653821914
$ws.Range("F21").Value = 36.19210460187975
$ws.Range("G21").Value = 3.608402624698795
$ws.Range("I21").Value = 52.25691677039583
$ws.Range("J21").Value = 9.867837687640529
$ws.Range("O21").Value = 26.91211612260957
$ws.Range("B22").Value = 19.58522931221873
$ws.Range("D22").Value = 9.032458281489077
$ws.Range("E22").Value = 13.70822871588955
$ws.Range("F22").Value = 36.58491920749592
$ws.Range("G22").Value = 3.605908622599556
$ws.Range("I22").Value = 53.9560509119212
$ws.Range("J22").Value = 9.876174675601838
$ws.Range("O22").Value = 27.17941864542292
$ws.Range("B23").Value = 19.27564066207565
$ws.Range("D23").Value = 9.032550767287978
$ws.Range("E23").Value = 13.68699813834026
$ws.Range("F23").Value = 36.37500870070366
$ws.Range("G23").Value = 3.607231046283911
$ws.Range("I23").Value = 53.05615300520156
$ws.Range("J23").Value = 9.871580611595775
$ws.Range("O23").Value = 27.0365083911838
$ws.Range("B24").Value = 18.04995794057753
$ws.Range("D24").Value = 9.035212148792816
$ws.Range("E24").Value = 13.61150965277894
$ws.Range("F24").Value = 35.58638200238035
$ws.Range("G24").Value = 3.612430596583006
$ws.Range("I24").Value = 49.49565175058866
$ws.Range("J24").Value = 9.85736025418872
$ws.Range("O24").Value = 26.50112243756142
$ws.Range("B25").Value = 16.61539234780995
$ws.Range("D25").Value = 9.043031288582753
$ws.Range("E25").Value = 13.54080724258617
$ws.Range("F25").Value = 34.75387368245401
$ws.Range("G25").Value = 3.618450967720129
$ws.Range("I25").Value = 45.33509062920486
$ws.Range("J25").Value = 9.848841359131361
$ws.Range("O25").Value = 25.93908695842768
